# games.xlsx — "Add files via upload"
#
# The sheet already had placeholder rows (season/week/date only) reserved
# through row 469; this upload fills in the week-13/14 game results that
# were missing (rows 454-481: team1/team2/score1/score2/home_team), and
# extends the sheet with additional still-empty placeholder rows (482-502)
# so the used range grows to A1:H502.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("games")

# Column C already uses a short-date style (numFmtId 14, style index 1) on
# every other row in the sheet. Copy that format onto the new date cells
# first so the cells we are about to fill reuse the same style instead of
# Excel minting a fresh (duplicate) style entry.
$ws.Cells.Item(454, 3).Copy() | Out-Null
$ws.Range("C467:C481").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# row -> week, date(serial), team1, team2, score1, score2
$games = @(
    @(454, 13, 45260, "SEA", "DAL", 35, 41),
    @(455, 13, 45263, "IND", "TEN", 31, 28),
    @(456, 13, 45263, "DET", "NO",  33, 28),
    @(457, 13, 45263, "ATL", "NYJ", 13, 8),
    @(458, 13, 45263, "LAC", "NE",  6,  0),
    @(459, 13, 45263, "ARI", "PIT", 24, 10),
    @(460, 13, 45263, "MIA", "WAS", 45, 15),
    @(461, 13, 45263, "DEN", "HOU", 17, 22),
    @(462, 13, 45263, "CAR", "TB",  18, 21),
    @(463, 13, 45263, "CLE", "LA",  19, 36),
    @(464, 13, 45263, "SF",  "PHI", 42, 19),
    @(465, 13, 45263, "KC",  "GB",  19, 27),
    @(466, 13, 45263, "CIN", "JAX", 34, 31),
    @(467, 14, 45267, "NE",  "PIT", 21, 18),
    @(468, 14, 45270, "CAR", "NO",  6,  28),
    @(469, 14, 45270, "DET", "CHI", 13, 28),
    @(470, 14, 45270, "IND", "CIN", 14, 34),
    @(471, 14, 45270, "LA",  "BAL", 31, 37),
    @(472, 14, 45270, "TB",  "ATL", 29, 25),
    @(473, 14, 45270, "HOU", "NYJ", 6,  30),
    @(474, 14, 45270, "JAX", "CLE", 27, 31),
    @(475, 14, 45270, "SEA", "SF",  16, 28),
    @(476, 14, 45270, "MIN", "LV",  3,  0),
    @(477, 14, 45270, "DEN", "LAC", 24, 7),
    @(478, 14, 45270, "BUF", "KC",  20, 17),
    @(479, 14, 45270, "PHI", "DAL", 13, 33),
    @(480, 14, 45271, "TEN", "MIA", 28, 27),
    @(481, 14, 45271, "GB",  "NYG", 22, 24)
)

foreach ($g in $games) {
    $r = $g[0]
    $ws.Cells.Item($r, 1).Value = 2023
    $ws.Cells.Item($r, 2).Value = $g[1]
    $ws.Cells.Item($r, 3).Value = $g[2]
    $ws.Cells.Item($r, 4).Value = $g[3]
    $ws.Cells.Item($r, 5).Value = $g[4]
    $ws.Cells.Item($r, 6).Value = $g[5]
    $ws.Cells.Item($r, 7).Value = $g[6]
    $ws.Cells.Item($r, 8).Formula = "=E$r"
}

# New trailing placeholder rows: season filled in, home_team formula
# carried down (columns B-G stay blank until future games are entered),
# matching rows 454-469 as they looked before this upload.
for ($r = 482; $r -le 501; $r++) {
    $ws.Cells.Item($r, 1).Value = 2023
    $ws.Cells.Item($r, 8).Formula = "=E$r"
}
$ws.Cells.Item(502, 1).Value = 2023

# Leave the workbook positioned/selected where the author last left it.
$ws.Range("B482").Select() | Out-Null
